# Apply "ADD gb.py and btes.py" change:
# For every year sheet (2025, 2030, 2035, 2040, 2045, 2050), a new technology
# column "gb" is inserted right after "eb" (becoming column B, pushing the
# rest right), and a new technology column "btes" is inserted right before
# "ites" (the last column). Row 2 investment-cost values are also refreshed
# to the new model results.

$wb = $excel.ActiveWorkbook

$sheetNames = @("2025", "2030", "2035", "2040", "2045", "2050")

# Final row-2 values (as strings, by NEW column letter) for each sheet,
# after both columns have been inserted.
$targets = @{
    "2025" = @{
        "A" = "0";                   "B" = "0";
        "C" = "348500.2913702198";   "D" = "0";
        "E" = "6952028.899878451";   "F" = "10001.38367865942";
        "G" = "0";                   "H" = "25342.77928792126";
        "I" = "0";                   "J" = "0";
        "K" = "0";                   "L" = "0";
        "M" = "0";                   "N" = "2171051.48368801";
        "O" = "19940.13531829329"
    }
    "2030" = @{
        "A" = "0";                   "B" = "0";
        "C" = "1146541.701244244";   "D" = "0";
        "E" = "0";                   "F" = "0";
        "G" = "0";                   "H" = "49137.49829535586";
        "I" = "0";                   "J" = "0";
        "K" = "0";                   "L" = "0";
        "M" = "0";                   "N" = "99597.19369803484";
        "O" = "37087.76506166223"
    }
    "2035" = @{
        "A" = "0";                   "B" = "0";
        "C" = "1041294.226359141";   "D" = "0";
        "E" = "0";                   "F" = "147293.0209330535";
        "G" = "0";                   "H" = "41476.88085549879";
        "I" = "0";                   "J" = "0";
        "K" = "0";                   "L" = "0";
        "M" = "0";                   "N" = "54681.56553335959";
        "O" = "39407.57948732926"
    }
    "2040" = @{
        "A" = "0";                   "B" = "0";
        "C" = "0";                   "D" = "0";
        "E" = "0";                   "F" = "0";
        "G" = "0";                   "H" = "0";
        "I" = "0";                   "J" = "0";
        "K" = "0";                   "L" = "0";
        "M" = "0";                   "N" = "3.645027391030453e-08";
        "O" = "0"
    }
    "2045" = @{
        "A" = "0";                   "B" = "0";
        "C" = "0";                   "D" = "0";
        "E" = "0";                   "F" = "0";
        "G" = "0";                   "H" = "0";
        "I" = "0";                   "J" = "0";
        "K" = "0";                   "L" = "0";
        "M" = "0";                   "N" = "0";
        "O" = "0"
    }
    "2050" = @{
        "A" = "0";                   "B" = "0";
        "C" = "0";                   "D" = "0";
        "E" = "0";                   "F" = "0";
        "G" = "0";                   "H" = "0";
        "I" = "0";                   "J" = "0";
        "K" = "0";                   "L" = "0";
        "M" = "0";                   "N" = "0";
        "O" = "0"
    }
}

# New (final) header row, left to right.
$headers = @("eb", "gb", "hp", "st", "wi", "ieh", "chp", "ac", "ab_ct", "ab_hp", "cp_ct", "cp_hp", "ttes", "btes", "ites")
$colLetters = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N", "O")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Insert "gb" as new column B (after "eb"), shifting hp..ites right by one.
    $ws.Columns.Item(2).Insert()

    # Insert "btes" as new column N (before "ites"), shifting ites right by one.
    # After the first insert, "ites" sits in column N (14th), so inserting at
    # column 14 pushes it to O and frees up N for "btes".
    $ws.Columns.Item(14).Insert()

    # Rewrite the full header row to the final arrangement.
    for ($i = 0; $i -lt $headers.Length; $i++) {
        $cellRef = $colLetters[$i] + "1"
        $ws.Range($cellRef).Value = $headers[$i]
    }

    # Rewrite row 2 values to the final (recomputed) investment-cost figures.
    $rowVals = $targets[$name]
    foreach ($col in $colLetters) {
        $cellRef = $col + "2"
        $ws.Range($cellRef).Value = [double]$rowVals[$col]
    }
}
